# Add a "LINCS protal URL" column (column R) to the Primary_Cells sheet,
# containing a link to the LINCS portal page for each row's pc_lincs_id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Primary_Cells_20171127")

# Header for the new column
$ws.Range("R1").Value = "LINCS protal URL"

# Find last used row in column A (data starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $lincsId = $ws.Cells.Item($r, 1).Text
    if ($lincsId) {
        $ws.Cells.Item($r, 18).Value = "http://lincsportal.ccs.miami.edu/cells/#/view/$lincsId"
    }
}

# Widen the new column to fit its (long URL) contents, like the author did.
$ws.Columns.Item(18).ColumnWidth = 46.3

# Scroll/selection state similar to author's saved view
$excel.ActiveWindow.ScrollColumn = 17
$ws.Range("U7").Select() | Out-Null
